$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = -21.9564
$ws.Range("A20").Value = -21.11509999999998
$ws.Range("A27").Value = -21.91689999999998
$ws.Range("A35").Value = -20.75359999999998
$ws.Range("A69").Value = -21.65899999999999
$ws.Range("A76").Value = -19.84939999999999
$ws.Range("A78").Value = -19.88829999999998
$ws.Range("A82").Value = -21.9246
$ws.Range("A83").Value = -21.806
$ws.Range("A93").Value = -21.03249999999999
